$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the time range text in B28 to reflect the finished session end time
$ws.Range("B28").Value = "14.00-16.00, 17.45-18.45, 19.00-20.00"

# Update the hours logged for that session (G28) from 3 to 4
$ws.Range("G28").Value = 4

# Recalculate so the SUM formula in H3 picks up the new total
$excel.Calculate()

# Reflect the user's final view/selection state: scrolled down and selected H28
$ws.Range("H28").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
